# Weekly update: insert two new "Berenjena" price records at the top of the
# existing Vega Central Mapocho de Santiago block (rows 157-181), shifting
# the older rows down by two (they become rows 159-183).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at row 157 (pushes old row 157.. down to 159..)
$ws.Rows.Item(157).Insert()
$ws.Rows.Item(157).Insert()

# New row 157: Berenjena, Primera, Región de Arica y Parinacota
$ws.Cells.Item(157, 1).Value = 9
$ws.Cells.Item(157, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(157, 3).Value = "Metropolitana"
$ws.Cells.Item(157, 4).Value = 44476
$ws.Cells.Item(157, 5).Value = 13
$ws.Cells.Item(157, 6).Value = 100112001
$ws.Cells.Item(157, 7).Value = "Berenjena"
$ws.Cells.Item(157, 8).Value = "Sin especificar"
$ws.Cells.Item(157, 9).Value = "Primera"
$ws.Cells.Item(157, 10).Value = 115
$ws.Cells.Item(157, 11).Value = 8000
$ws.Cells.Item(157, 12).Value = 10000
$ws.Cells.Item(157, 13).Value = 9009
$ws.Cells.Item(157, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(157, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(157, 16).Value = 150
$ws.Cells.Item(157, 17).Value = 60
$ws.Cells.Item(157, 18).Value = "Hortaliza"

# New row 158: Berenjena, Segunda, Región de Arica y Parinacota
$ws.Cells.Item(158, 1).Value = 9
$ws.Cells.Item(158, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(158, 3).Value = "Metropolitana"
$ws.Cells.Item(158, 4).Value = 44476
$ws.Cells.Item(158, 5).Value = 13
$ws.Cells.Item(158, 6).Value = 100112001
$ws.Cells.Item(158, 7).Value = "Berenjena"
$ws.Cells.Item(158, 8).Value = "Sin especificar"
$ws.Cells.Item(158, 9).Value = "Segunda"
$ws.Cells.Item(158, 10).Value = 43
$ws.Cells.Item(158, 11).Value = 7000
$ws.Cells.Item(158, 12).Value = 7000
$ws.Cells.Item(158, 13).Value = 7000
$ws.Cells.Item(158, 14).Value = "$/caja 100 unidades"
$ws.Cells.Item(158, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(158, 16).Value = 70
$ws.Cells.Item(158, 17).Value = 100
$ws.Cells.Item(158, 18).Value = "Hortaliza"
